# Update activity log: add new row 7 with the latest entry,
# and select the newly added cell D7 (mirrors previous D6 selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data (A7, C7, D7) - B7 stays blank, matching row 6's pattern.
$ws.Range("A7").Value = "2/13/2020jaclemon"
$ws.Range("C7").Value = "20 minutes"
$ws.Range("D7").Value = "Used Clion to fix command line issue where -records would read as -r"

# Update selection to the new last cell, like Excel does after data entry.
$ws.Range("D7").Select()
